$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step A: insert a new row before the current last row (row 13: 4x5090, 2501.38...)
# so that row 13 becomes the new "4x5090, 9790.98" entry and the old row13 shifts to row 14.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "4x5090"
$ws.Range("B13").Value = 9790.98
$ws.Range("C13").Value = 2.6
$ws.Range("D13").Value = 0.07376403814758301

# Step B: append a new row after the (now shifted) old last row at row 14, as row 15.
$ws.Range("A15").Value = "4x5090"
$ws.Range("B15").Value = 4622.08
$ws.Range("C15").Value = 2.6
$ws.Range("D15").Value = 0.1562548078402412

# Step C: insert a new row before row 6 (currently "2x4090, 467.26"), pushing
# rows 6-15 down to 7-16, and fill it with the new "1x6000" entry.
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "1x6000"
$ws.Range("B6").Value = 3168.82
$ws.Range("C6").Value = 1.29
$ws.Range("D6").Value = 0.1130809996570753
